$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price/Volume columns store plain numeric-looking text (e.g. "331.13", "7.51%").
# Force NumberFormat to Text ("@") before assignment so Excel keeps them as
# literal strings instead of auto-converting to numbers/percentages.

$ws.Range("D2:E2").NumberFormat = "@"
$ws.Range("D2").Value = "331.13"
$ws.Range("E2").Value = "7.51%"

$ws.Range("D3:E3").NumberFormat = "@"
$ws.Range("D3").Value = "40.73"
$ws.Range("E3").Value = "12.83%"

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "17.58%"

$ws.Range("D5:E5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08208"
$ws.Range("E5").Value = "6.57%"

$ws.Range("D6:E6").NumberFormat = "@"
$ws.Range("D6").Value = "8.804"
$ws.Range("E6").Value = "6.20%"

$ws.Range("D7:E7").NumberFormat = "@"
$ws.Range("D7").Value = "4.574"
$ws.Range("E7").Value = "4.34%"

$ws.Range("D8:E8").NumberFormat = "@"
$ws.Range("D8").Value = "1.984"
$ws.Range("E8").Value = "7.51%"

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.10%"

$ws.Range("D10:E10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9505"
$ws.Range("E10").Value = "3.25%"

$ws.Range("D11:E11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1351"
$ws.Range("E11").Value = "23.11%"

$ws.Range("D12:E12").NumberFormat = "@"
$ws.Range("D12").Value = "0.2000"
$ws.Range("E12").Value = "8.52%"

$ws.Range("D13:E13").NumberFormat = "@"
$ws.Range("B13").Value = "MCDex"
$ws.Range("C13").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D13").Value = "10.55"
$ws.Range("E13").Value = "66.66%"

$ws.Range("D14:E14").NumberFormat = "@"
$ws.Range("B14").Value = "MandalaExchangeToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D14").Value = "0.09252"
$ws.Range("E14").Value = "5.57%"

$ws.Range("D15:E15").NumberFormat = "@"
$ws.Range("B15").Value = "BitrueCoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D15").Value = "0.03498"
$ws.Range("E15").Value = "4.57%"

$ws.Range("D16:E16").NumberFormat = "@"
$ws.Range("B16").Value = "BitMartToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D16").Value = "0.09635"
$ws.Range("E16").Value = "1.10%"

$ws.Range("D17:E17").NumberFormat = "@"
$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D17").Value = "0.001305"
$ws.Range("E17").Value = "-5.33%"

$ws.Range("D18:E18").NumberFormat = "@"
$ws.Range("D18").Value = "0.006261"
$ws.Range("E18").Value = "1.02%"

$ws.Range("D19:E19").NumberFormat = "@"
$ws.Range("D19").Value = "3.352"
$ws.Range("E19").Value = "-0.29%"

$ws.Range("D20:E20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3538"
$ws.Range("E20").Value = "2.99%"

$ws.Range("D21:E21").NumberFormat = "@"
$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D21").Value = "0.1445"
$ws.Range("E21").Value = "12.05%"

$ws.Range("D22:E22").NumberFormat = "@"
$ws.Range("B22").Value = "ZBToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D22").Value = "0.2432"
$ws.Range("E22").Value = "5.37%"

$ws.Range("D23:E23").NumberFormat = "@"
$ws.Range("B23").Value = "CoinExToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D23").Value = "0.04435"
$ws.Range("E23").Value = "2.70%"

$ws.Range("D24:E24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001254"
$ws.Range("E24").Value = "4.66%"

$ws.Range("D25:E25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004460"
$ws.Range("E25").Value = "4.71%"

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-10.56%"

$ws.Range("D27:E27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003977"
$ws.Range("E27").Value = "37.12%"

$ws.Range("D39:E39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02503"
$ws.Range("E39").Value = "20.12%"

$ws.Range("D40:E40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05298"
$ws.Range("E40").Value = "6.80%"

$ws.Range("D41:E41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007462"
$ws.Range("E41").Value = "-0.60%"

$ws.Range("D42:E42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1453"
$ws.Range("E42").Value = "7.64%"

$ws.Range("D43:E43").NumberFormat = "@"
$ws.Range("D43").Value = "0.008961"
$ws.Range("E43").Value = "6.80%"

$ws.Range("D44:E44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002045"
$ws.Range("E44").Value = "-1.24%"

$ws.Range("D45:E45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01059"
$ws.Range("E45").Value = "26.13%"

$ws.Range("D46:E46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006785"
$ws.Range("E46").Value = "7.63%"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.37%"

$ws.Range("D48:E48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003464"
$ws.Range("E48").Value = "21.51%"

$ws.Range("D49:E49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001794"
$ws.Range("E49").Value = "6.20%"

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.37%"

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.37%"
